$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "variable" field renamed to "key" throughout the metadata sheet:
#  - header cell C1: "variable" -> "key"
#  - template cells E2, E3: "{variable}" -> "{key}"
$ws.Range("C1").Value = "key"
$ws.Range("E2").Value = "{key}"
$ws.Range("E3").Value = "{key}"

# selection moves to E4
$ws.Range("E4").Select()
